$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "Contratto"
$ws.Range("E3").Value = "Contratto"

$ws.Range("M3").Select()
